# Auto-generated script to apply Chocobo_Profits.xlsx numeric updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 52
$ws.Range("H52").Value = 102500
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 102500
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 307500
$ws.Range("N52").Value = -307820
$ws.Range("M52").ClearContents()
# Row 58
$ws.Range("H58").Value = 3543.3
$ws.Range("I58").Value = 603.6667
$ws.Range("K58").Value = 1811.0001
$ws.Range("M58").Value = -1661.0001
# Row 98
$ws.Range("H98").Value = 9308.714
$ws.Range("I98").Value = 8603.637000000001
$ws.Range("J98").Value = 9764.941000000001
$ws.Range("K98").Value = 8603.637000000001
$ws.Range("L98").Value = 9764.941000000001
$ws.Range("M98").Value = -7105.637000000001
$ws.Range("N98").Value = -12760.941
# Row 122
$ws.Range("H122").Value = 9308.714
$ws.Range("I122").Value = 8603.637000000001
$ws.Range("J122").Value = 9764.941000000001
$ws.Range("K122").Value = 25810.911
$ws.Range("L122").Value = 29294.823
$ws.Range("M122").Value = -23360.911
$ws.Range("N122").Value = -34194.823
# Row 129
$ws.Range("H129").Value = 860.96
$ws.Range("J129").Value = 864.6061
$ws.Range("L129").Value = 2593.8183
$ws.Range("N129").Value = -12593.8183
# Row 138
$ws.Range("H138").Value = 2738.4895
$ws.Range("I138").Value = 1466
$ws.Range("J138").Value = 2992.9875
$ws.Range("K138").Value = 4398
$ws.Range("L138").Value = 8978.962500000001
$ws.Range("M138").Value = 742
$ws.Range("N138").Value = -19258.9625

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6682.4263
$ws.Range("I32").Value = 5580.343
$ws.Range("J32").Value = 8166
$ws.Range("K32").Value = 5580.343
$ws.Range("L32").Value = 8166
$ws.Range("M32").Value = -5293.343
$ws.Range("N32").Value = -8740
# Row 92
$ws.Range("H92").Value = 14366.667
$ws.Range("J92").Value = 14366.667
$ws.Range("L92").Value = 14366.667
$ws.Range("N92").Value = -19358.667
# Row 108
$ws.Range("H108").Value = 38684
$ws.Range("J108").Value = 38684
$ws.Range("L108").Value = 38684
$ws.Range("N108").Value = -46364
# Row 115
$ws.Range("H115").Value = 29950
$ws.Range("J115").Value = 29950
$ws.Range("L115").Value = 29950
$ws.Range("N115").Value = -33084
# Row 132
$ws.Range("H132").Value = 2307.4211
$ws.Range("I132").Value = 1024.5
$ws.Range("J132").Value = 5899.6
$ws.Range("K132").Value = 3073.5
$ws.Range("L132").Value = 17698.8
$ws.Range("M132").Value = -543.5
$ws.Range("N132").Value = -22758.8

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2423.5557
$ws.Range("I86").Value = 2218.6667
$ws.Range("J86").Value = 2833.3333
$ws.Range("K86").Value = 2218.6667
$ws.Range("L86").Value = 2833.3333
$ws.Range("M86").Value = -1095.6667
$ws.Range("N86").Value = -5079.3333
# Row 89
$ws.Range("H89").Value = 2423.5557
$ws.Range("I89").Value = 2218.6667
$ws.Range("J89").Value = 2833.3333
$ws.Range("K89").Value = 11093.3335
$ws.Range("L89").Value = 14166.6665
$ws.Range("M89").Value = -5477.333500000001
$ws.Range("N89").Value = -25398.6665
# Row 95
$ws.Range("H95").Value = 32600
$ws.Range("J95").Value = 32600
$ws.Range("L95").Value = 32600
$ws.Range("N95").Value = -38092

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3911.45
$ws.Range("I31").Value = 1436.74
$ws.Range("J31").Value = 16285
$ws.Range("K31").Value = 1436.74
$ws.Range("L31").Value = 16285
$ws.Range("M31").Value = -1141.74
$ws.Range("N31").Value = -16875
# Row 34
$ws.Range("H34").Value = 3911.45
$ws.Range("I34").Value = 1436.74
$ws.Range("J34").Value = 16285
$ws.Range("K34").Value = 1436.74
$ws.Range("L34").Value = 16285
$ws.Range("M34").Value = -1234.74
$ws.Range("N34").Value = -16689
# Row 58
$ws.Range("H58").Value = 1934.3334
$ws.Range("I58").Value = 1465.3572
$ws.Range("J58").Value = 8500
$ws.Range("K58").Value = 1465.3572
$ws.Range("L58").Value = 8500
$ws.Range("M58").Value = -1262.3572
$ws.Range("N58").Value = -8906
# Row 136
$ws.Range("H136").Value = 1934.3334
$ws.Range("I136").Value = 1465.3572
$ws.Range("J136").Value = 8500
$ws.Range("K136").Value = 4396.071599999999
$ws.Range("L136").Value = 25500
$ws.Range("M136").Value = -1846.071599999999
$ws.Range("N136").Value = -30600
# Row 137
$ws.Range("H137").Value = 45320
$ws.Range("J137").Value = 45320
$ws.Range("L137").Value = 45320
$ws.Range("N137").Value = -55520

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 1280.2142
$ws.Range("I22").Value = 453.8
$ws.Range("J22").Value = 1739.3334
$ws.Range("K22").Value = 1361.4
$ws.Range("L22").Value = 5218.0002
$ws.Range("M22").Value = -1192.4
$ws.Range("N22").Value = -5556.0002
# Row 27
$ws.Range("H27").Value = 1280.2142
$ws.Range("I27").Value = 453.8
$ws.Range("J27").Value = 1739.3334
$ws.Range("K27").Value = 1361.4
$ws.Range("L27").Value = 5218.0002
$ws.Range("M27").Value = -1259.4
$ws.Range("N27").Value = -5422.0002
# Row 48
$ws.Range("H48").Value = 8250
$ws.Range("J48").Value = 8250
$ws.Range("L48").Value = 24750
$ws.Range("N48").Value = -25250
# Row 55
$ws.Range("H55").Value = 4791.7646
$ws.Range("J55").Value = 5028.75
$ws.Range("L55").Value = 15086.25
$ws.Range("N55").Value = -15440.25
# Row 58
$ws.Range("H58").Value = 3348
$ws.Range("J58").Value = 3997.5
$ws.Range("L58").Value = 11992.5
$ws.Range("N58").Value = -12248.5
# Row 131
$ws.Range("H131").Value = 782.24
$ws.Range("I131").Value = 305
$ws.Range("J131").Value = 823.73914
$ws.Range("K131").Value = 915
$ws.Range("L131").Value = 2471.21742
$ws.Range("M131").Value = 4125
$ws.Range("N131").Value = -12551.21742

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 10468792
$ws.Range("J11").Value = 781300.3
$ws.Range("L11").Value = 781300.3
$ws.Range("N11").Value = -781578.3
# Row 132
$ws.Range("H132").Value = 6287.143
$ws.Range("I132").Value = 3670.6667
$ws.Range("K132").Value = 11012.0001
$ws.Range("M132").Value = -8482.000100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 32613.625
$ws.Range("I22").Value = 48830.57
$ws.Range("J22").Value = 1654
$ws.Range("K22").Value = 48830.57
$ws.Range("L22").Value = 1654
$ws.Range("M22").Value = -48535.57
$ws.Range("N22").Value = -2244
# Row 27
$ws.Range("H27").Value = 32613.625
$ws.Range("I27").Value = 48830.57
$ws.Range("J27").Value = 1654
$ws.Range("K27").Value = 48830.57
$ws.Range("L27").Value = 1654
$ws.Range("M27").Value = -48723.57
$ws.Range("N27").Value = -1868
# Row 61
$ws.Range("H61").Value = 1385.5714
$ws.Range("I61").Value = 1268.625
$ws.Range("J61").Value = 1759.8
$ws.Range("K61").Value = 1268.625
$ws.Range("L61").Value = 1759.8
$ws.Range("M61").Value = -1066.625
$ws.Range("N61").Value = -2163.8
# Row 113
$ws.Range("H113").Value = 1385.5714
$ws.Range("I113").Value = 1268.625
$ws.Range("J113").Value = 1759.8
$ws.Range("K113").Value = 1268.625
$ws.Range("L113").Value = 1759.8
$ws.Range("M113").Value = 901.375
$ws.Range("N113").Value = -6099.8
# Row 132
$ws.Range("H132").Value = 6503.16
$ws.Range("I132").Value = 2425.4
$ws.Range("K132").Value = 7276.200000000001
$ws.Range("M132").Value = -4746.200000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 52245.445
$ws.Range("J46").Value = 52245.445
$ws.Range("L46").Value = 52245.445
$ws.Range("N46").Value = -52707.445
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# Row 132
$ws.Range("H132").Value = 15880025
$ws.Range("I132").Value = 9339.416999999999
$ws.Range("K132").Value = 28018.251
$ws.Range("M132").Value = -25488.251
# Row 134
$ws.Range("H134").Value = 52245.445
$ws.Range("J134").Value = 52245.445
$ws.Range("L134").Value = 156736.335
$ws.Range("N134").Value = -161806.335

